# Apply updated dSF (column F) values for the specified rows.
# Mapping of worksheet row number -> new value in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 4
    11 = 1
    12 = 3
    13 = -2
    16 = -7
    19 = -3
    20 = 1
    21 = 2
    22 = 2
    24 = -3
    27 = 4
    29 = 2
    31 = -2
    32 = 0
    34 = 0
    35 = -2
    36 = -1
    39 = -4
    43 = -1
    45 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
